$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9535946846008301
$ws.Range("B1").Value = 2.085201025009155
$ws.Range("C1").Value = 7.968392372131348
$ws.Range("D1").Value = 2.522231101989746
$ws.Range("E1").Value = 0.679522693157196
